# Update "run 62" optimisation results: Schedule (pump run summary) and Detailed (price/type/pump_status) sheets
$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- "Schedule" sheet: recomputed pump-run windows (rows 2-4 replaced with latest optimisation, former row 5 removed) ---
$wsSchedule.Range("A2").Value = 46039.29166666666
$wsSchedule.Range("B2").Value = 46039.47916666666
$wsSchedule.Range("C2").Value = 4.5
$wsSchedule.Range("D2").Value = 17.01
$wsSchedule.Range("E2").Value = 278.8351215
$wsSchedule.Range("F2").Value = 16.3924233686067

$wsSchedule.Range("A3").Value = 46039.52083333334
$wsSchedule.Range("B3").Value = 46040
$wsSchedule.Range("C3").Value = 11.5
$wsSchedule.Range("D3").Value = 43.47
$wsSchedule.Range("E3").Value = 537.543513
$wsSchedule.Range("F3").Value = 12.36585031055901

$wsSchedule.Range("A4").Value = 46040.29166666666
$wsSchedule.Range("B4").Value = 46040.79166666666
$wsSchedule.Range("C4").Value = 12
$wsSchedule.Range("D4").Value = 45.36
$wsSchedule.Range("E4").Value = 25.08997725
$wsSchedule.Range("F4").Value = 0.553130009920635

# The optimisation horizon shrank by one scheduled run, so the last row is no longer present
$wsSchedule.Rows.Item(5).Delete()

# --- "Detailed" sheet: updated Price values, a couple of Type flips, and Pump_Status toggles ---
$wsDetailed.Range("E3").Value = "OFF"
$wsDetailed.Range("E4").Value = "OFF"
$wsDetailed.Range("E5").Value = "OFF"
$wsDetailed.Range("E6").Value = "OFF"
$wsDetailed.Range("E7").Value = "OFF"
$wsDetailed.Range("E8").Value = "OFF"
$wsDetailed.Range("E9").Value = "OFF"
$wsDetailed.Range("E10").Value = "OFF"
$wsDetailed.Range("E25").Value = "OFF"
$wsDetailed.Range("B33").Value = 36.06038
$wsDetailed.Range("B34").Value = 36.06036
$wsDetailed.Range("B35").Value = -7.37172
$wsDetailed.Range("C35").Value = "historical"
$wsDetailed.Range("B36").Value = -6.48811
$wsDetailed.Range("C36").Value = "historical"
$wsDetailed.Range("B37").Value = -2.93457
$wsDetailed.Range("B38").Value = -2.81393
$wsDetailed.Range("B39").Value = -2.96975
$wsDetailed.Range("B40").Value = 2.29449
$wsDetailed.Range("B41").Value = 13.64807
$wsDetailed.Range("E41").Value = "ON"
$wsDetailed.Range("B42").Value = 33.18906
$wsDetailed.Range("E42").Value = "ON"
$wsDetailed.Range("B43").Value = 36.2
$wsDetailed.Range("E43").Value = "ON"
$wsDetailed.Range("B44").Value = 30.11497
$wsDetailed.Range("E44").Value = "ON"
$wsDetailed.Range("B45").Value = 36.0601
$wsDetailed.Range("E45").Value = "ON"
$wsDetailed.Range("B46").Value = 36.05843
$wsDetailed.Range("E46").Value = "ON"
$wsDetailed.Range("B47").Value = 36.06041
$wsDetailed.Range("E47").Value = "ON"
$wsDetailed.Range("B48").Value = 36.06043
$wsDetailed.Range("E48").Value = "ON"
$wsDetailed.Range("B49").Value = 36.0604
$wsDetailed.Range("E49").Value = "ON"
$wsDetailed.Range("B50").Value = 36.06038
$wsDetailed.Range("B52").Value = 36.0603
$wsDetailed.Range("B53").Value = 36.06028
$wsDetailed.Range("B54").Value = 36.06027
$wsDetailed.Range("B55").Value = 36.06026
$wsDetailed.Range("B56").Value = 56.97996
$wsDetailed.Range("B57").Value = 56.98
$wsDetailed.Range("B58").Value = 56.98
$wsDetailed.Range("B59").Value = 56.98
$wsDetailed.Range("B62").Value = 56.98
$wsDetailed.Range("B64").Value = 36.0595
$wsDetailed.Range("B65").Value = 23.49478
$wsDetailed.Range("B66").Value = -0.90707
$wsDetailed.Range("B67").Value = 0.7
$wsDetailed.Range("B68").Value = 0.7
$wsDetailed.Range("B69").Value = 0.51
$wsDetailed.Range("B70").Value = 0.51
$wsDetailed.Range("B71").Value = 35.88
$wsDetailed.Range("B72").Value = 0.66409
$wsDetailed.Range("B73").Value = -0.97534
$wsDetailed.Range("B74").Value = 0.0003
$wsDetailed.Range("B75").Value = -0.92843
$wsDetailed.Range("B76").Value = -5.51
$wsDetailed.Range("B77").Value = -5.51
$wsDetailed.Range("B78").Value = -6.8
$wsDetailed.Range("B79").Value = -7
$wsDetailed.Range("B80").Value = -8.06785
$wsDetailed.Range("B81").Value = -6.77128
$wsDetailed.Range("B82").Value = -7.21403
$wsDetailed.Range("B83").Value = -7.00224
$wsDetailed.Range("B84").Value = -6.52915
$wsDetailed.Range("B85").Value = -6.51864
$wsDetailed.Range("B86").Value = -3.05165
$wsDetailed.Range("B87").Value = 0.00032
$wsDetailed.Range("B89").Value = 46.21963
$wsDetailed.Range("B90").Value = 54.93629
$wsDetailed.Range("B91").Value = 55.52267
$wsDetailed.Range("B92").Value = 46.9203
$wsDetailed.Range("B94").Value = 47.36052
